# Add two new columns (G: "Expenses type", H: "Saving opportunity") that
# mirror / derive from the existing "Expenses Type" (E) and
# "Savings Oppurtunity" (F) columns, for the monthly expenses sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- Header row ------------------------------------------------------
$ws.Range("G1").Value = "Expenses type"
$ws.Range("H1").Value = "Saving opportunity"

# Match the header formatting already used by the other header cells
# (bold font, thin border, centered) by copying F1's format onto the
# two new header cells.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows ---------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $expensesType = $ws.Cells.Item($r, 5).Value2      # column E
    $savingsOpp   = $ws.Cells.Item($r, 6).Value2       # column F

    # G mirrors the existing "Expenses Type" value.
    $ws.Cells.Item($r, 7).Value = $expensesType

    # H is "Yes" only when the expense is High-type AND already flagged
    # as a savings opportunity; otherwise "No".
    if ($expensesType -eq "High" -and $savingsOpp -eq "Yes") {
        $ws.Cells.Item($r, 8).Value = "Yes"
    } else {
        $ws.Cells.Item($r, 8).Value = "No"
    }
}

